# refreshed logistic reg with kfolds
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the RAW row (row 4) accuracy values with the new k-folds results ---
# Order matters for shared-string table layout: set F4 before C4.
$ws.Range("F4").Value = "0.652`nfalse pos. 44, false n. 0"
$ws.Range("C4").Value = "0.903`nfalse pos. 2 false n. 10"

# Wrap the long confusion-matrix text in the RAW row
$ws.Range("C4").WrapText = $true
$ws.Range("F4").WrapText = $true

# --- Row heights: rows 3-8 grow to fit the new wrapped text ---
$ws.Rows.Item(3).RowHeight = 44.25
$ws.Rows.Item(4).RowHeight = 44.25
$ws.Rows.Item(5).RowHeight = 44.25
$ws.Rows.Item(6).RowHeight = 44.25
$ws.Rows.Item(7).RowHeight = 44.25
$ws.Rows.Item(8).RowHeight = 44.25

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666
$ws.Columns.Item(3).ColumnWidth = 23.5
$ws.Columns.Item(4).ColumnWidth = 23.5
$ws.Columns.Item(5).ColumnWidth = 23.5
$ws.Columns.Item(6).ColumnWidth = 23.5

# --- New scratch cell below the table, formatted as text ---
$ws.Range("C12").NumberFormat = "@"

# --- Selection moves to C14 ---
$null = $ws.Range("C14").Select()
